# Header text updates: bulletin volume/number and week-covering dates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# Crime-complaints table updates (rows 14-31)
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("M14").Copy($ws.Range("E14"))
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("L15").Value = 0
$ws.Range("L15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N15").Value = 0
$ws.Range("N15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 37.5
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 14.285714285714
$ws.Range("L16").Value = 60
$ws.Range("M16").Value = -46.666666666666
$ws.Range("N16").Value = -85.714285714285
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 400
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = -30
$ws.Range("L17").Value = -61.111111111111
$ws.Range("M17").Value = 75
$ws.Range("N17").Value = -73.076923076923
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 76.923076923076
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 75
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = 600
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N18").Value = -61.111111111111
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 9.090909090909
$ws.Range("I19").Value = 33
$ws.Range("J19").Value = 29
$ws.Range("K19").Value = 13.793103448275
$ws.Range("L19").Value = 26.923076923076
$ws.Range("M19").Value = 57.142857142857
$ws.Range("N19").Value = -29.787234042553
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 4
$ws.Range("K20").Value = 300
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -89.189189189189
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 44.444444444444
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = 12.222222222222
$ws.Range("I21").Value = 67
$ws.Range("J21").Value = 56
$ws.Range("K21").Value = 19.642857142857
$ws.Range("L21").Value = 6.349206349206
$ws.Range("M21").Value = 48.888888888888
$ws.Range("N21").Value = -66.995073891625
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 2
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -75
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = 100
$ws.Range("M23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 22.222222222222
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = -11.111111111111
$ws.Range("I24").Value = 71
$ws.Range("J24").Value = 69
$ws.Range("K24").Value = 2.898550724637
$ws.Range("L24").Value = -20.224719101123
$ws.Range("M24").Value = -25.263157894736
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 73
$ws.Range("G25").Value = 98
$ws.Range("H25").Value = -25.510204081632
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 59
$ws.Range("K25").Value = -18.64406779661
$ws.Range("L25").Value = -42.168674698795
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 40
$ws.Range("I26").Value = 23
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = 27.777777777777
$ws.Range("L26").Value = -11.538461538461
$ws.Range("M26").Value = 9.523809523809
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 2
$ws.Range("L27").Value = 100
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("L28").Value = -33.333333333333
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("M14").Copy($ws.Range("E29"))
$ws.Range("L29").Value = -100
$ws.Range("L29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("M14").Copy($ws.Range("E30"))
$ws.Range("L30").Value = -100
$ws.Range("L30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 2
$ws.Range("L31").Value = -100
$ws.Range("L31").NumberFormat = "#,##0.0;""-""#,##0.0"
